# "Offline mode" feature:
# Populate the degrees/minutes/seconds breakdown columns (B:G) for every
# coordinate row, and refresh the EGSA87 (X,Y) projected-coordinate columns
# (J:K) so the sheet can be used without recomputing the transform online.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Degrees / minutes / seconds breakdown for latitude (B,C,D) and longitude (E,F,G)
# for rows 3-13, matching the decimal values already present in H:I.
$dms = @{
    3  = @(40,31,24.49, 21,41,4.34)
    4  = @(40,31,24.21, 21,41,4.38)
    5  = @(40,31,24.19, 21,41,4.15)
    6  = @(40,31,25.85, 21,40,59.87)
    7  = @(40,31,25.73, 21,40,59.84)
    8  = @(40,31,25.72, 21,40,59.71)
    9  = @(40,31,25.63, 21,40,59.71)
    10 = @(40,31,25.52, 21,40,59.74)
    11 = @(40,31,25.5,  21,40,59.7)
    12 = @(40,31,24.27, 21,40,59.45)
    13 = @(40,31,24.27, 21,40,59.36)
}

foreach ($r in 3..13) {
    $vals = $dms[$r]
    $ws.Range("B$r").Value = $vals[0]
    $ws.Range("C$r").Value = $vals[1]
    $ws.Range("D$r").Value = $vals[2]
    $ws.Range("E$r").Value = $vals[3]
    $ws.Range("F$r").Value = $vals[4]
    $ws.Range("G$r").Value = $vals[5]
}

# Refresh the EGSA87 projected coordinates (J,K) for rows 6-13 offline.
$jk = @{
    6  = @("303707.14658056", "4488140.5979287")
    7  = @("303607.77423981", "4488194.4349555")
    8  = @("303607.77423981", "4488194.4349555")
    9  = @("303606.97129996", "4488190.7533626")
    10 = @("303606.97129996", "4488190.7533626")
    11 = @("303606.97129996", "4488190.7533626")
    12 = @("303606.97129996", "4488190.7533626")
    13 = @("303603.90411498", "4488190.5251506")
}

# These columns hold their EGSA87 figures as text (same as the surrounding
# H:K cells already on the sheet), so enter them with a leading apostrophe to
# keep them as text instead of letting Excel coerce them to numbers.
foreach ($r in 6..13) {
    $pair = $jk[$r]
    $ws.Range("J$r").Value = "'" + $pair[0]
    $ws.Range("K$r").Value = "'" + $pair[1]
}

# Update the selected range to reflect the newly populated columns.
[void]$ws.Range("H3:K13").Select()
